# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Vega Modelo de Temuco - Espárragos"
# at the top of the data block (rows 83-84), pushing the existing rows
# 83-94 down to rows 85-96.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before the current row 83, shifting rows 83:94 down to 85:96.
$ws.Range("A83:R84").EntireRow.Insert()

# --- New row 83 ---
$ws.Range("A83").Value = 10
$ws.Range("B83").Value = "Vega Modelo de Temuco"
$ws.Range("C83").Value = "La Araucanía"
$ws.Range("D83").Value2 = 45211
$ws.Range("E83").Value = 9
$ws.Range("F83").Value = 300000000
$ws.Range("G83").Value = "Espárragos"
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 400
$ws.Range("K83").Value = 1500
$ws.Range("L83").Value = 1500
$ws.Range("M83").Value = 1500
$ws.Range("N83").Value = "$/kilo"
$ws.Range("O83").Value = "Región de La Araucanía"
$ws.Range("P83").Value = 1500
$ws.Range("Q83").Value = 1
$ws.Range("R83").Value = "Hortaliza"

# --- New row 84 ---
$ws.Range("A84").Value = 10
$ws.Range("B84").Value = "Vega Modelo de Temuco"
$ws.Range("C84").Value = "La Araucanía"
$ws.Range("D84").Value2 = 45211
$ws.Range("E84").Value = 9
$ws.Range("F84").Value = 300000000
$ws.Range("G84").Value = "Espárragos"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 500
$ws.Range("K84").Value = 1600
$ws.Range("L84").Value = 1600
$ws.Range("M84").Value = 1600
$ws.Range("N84").Value = "$/kilo"
$ws.Range("O84").Value = "Región del Maule"
$ws.Range("P84").Value = 1600
$ws.Range("Q84").Value = 1
$ws.Range("R84").Value = "Hortaliza"
